$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.651.34'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.348.20'
$ws.Range("E3").Value = '  -1.78%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.34'
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '659.53'
$ws.Range("E6").Value = '  +0.58%  '

$ws.Range("E7").Value = '  -3.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.425'
$ws.Range("E8").Value = '  -2.52%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("E10").Value = '  -5.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.346.52'
$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("E12").Value = '  -2.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.49'
$ws.Range("E13").Value = '  -0.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '97.458.36'
$ws.Range("E14").Value = '  +0.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.11'
$ws.Range("E15").Value = '  -5.13%  '

$ws.Range("E16").Value = '  -2.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.969.15'
$ws.Range("E17").Value = '  -1.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.87'
$ws.Range("E18").Value = '  +2.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.399.32'
$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.85'
$ws.Range("E20").Value = '  +1.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.555'
$ws.Range("E21").Value = '  +10.45%  '

$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '511.71'
$ws.Range("E23").Value = '  +0.72%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.35'
$ws.Range("E24").Value = '  -3.25%  '

$ws.Range("E25").Value = '  -3.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.63'
$ws.Range("E26").Value = '  +6.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '97.50'
$ws.Range("E27").Value = '  -1.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.27'
$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("E29").Value = '  -3.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.66'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("E32").Value = '  -5.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.56'
$ws.Range("E33").Value = '  +11.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.566'
$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.65'
$ws.Range("E36").Value = '  -3.82%  '

$ws.Range("E37").Value = '  +2.58%  '

$ws.Range("E38").Value = '  +4.96%  '

$ws.Range("E39").Value = '  -0.26%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '515.63'
$ws.Range("E41").Value = '  -3.13%  '

$ws.Range("E42").Value = '  -1.25%  '

$ws.Range("E43").Value = '  +3.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.847'
$ws.Range("E44").Value = '  -1.89%  '

$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.72'
$ws.Range("E45").Value = '  +7.76%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.71'
$ws.Range("E46").Value = '  +3.57%  '

$ws.Range("E47").Value = '  +4.37%  '

$ws.Range("E48").Value = '  -1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.86'
$ws.Range("E49").Value = '  +7.02%  '

$ws.Range("E50").Value = '  -4.84%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '164.06'
$ws.Range("E51").Value = '  +1.31%  '
